# Update the "Module 4" agenda line on the closing slide:
#   "Module 4`tUse Visual Studio Code Git integration "
# becomes
#   "Module 4`tBasics of Visual Studio Code Git integration "
#
# The run also contains a second, separately-formatted run holding the
# "(15 min)" duration which must stay untouched, so we locate the exact
# sub-range of characters belonging to the first run (via its text) and
# replace just that sub-range's text instead of overwriting the whole
# shape / paragraph (which would collapse formatting and other runs).

$p = $ppt.ActivePresentation

$oldText = "Module 4`tUse Visual Studio Code Git integration "
$newText = "Module 4`tBasics of Visual Studio Code Git integration "

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if (-not $sh.HasTextFrame) { continue }
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf($oldText)
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, $oldText.Length)
            $sub.Text = $newText
        }
    }
}
